$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally listed, for each terminator, three rows: the base
# sequence, a "+ 250bp" variant, and a "+ 500bp" variant. Keep only the base
# terminator rows (SV40, bGH225, sT1A18, T1A9) and drop the extended variants,
# leaving a simple 2-column reference table (Terminator / Sequence).
$rowsToDelete = @(13, 12, 10, 9, 7, 6, 4, 3)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
